# Rework the merge-field table row so that the {#items} loop-open marker
# only appears in the first ("description") cell and the {/items}
# loop-close marker only appears in the last ("finalAmount") cell.
#
#   {#items}{description}{/items}  -> {#items}{description}
#   {#items}{HSN}{/items}          -> {HSN}
#   {#items}{GST}{/items}          -> {GST}
#   {#items}{quantity}{/items}     -> {quantity}
#   {#items}{price}{/items}        -> {price}
#   {#items}{finalAmount}{/items}  -> {finalAmount}{/items}

$d = $word.ActiveDocument

$d.Content.Find.Execute("{#items}{description}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{#items}{description}", 2)
$d.Content.Find.Execute("{#items}{HSN}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{HSN}", 2)
$d.Content.Find.Execute("{#items}{GST}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{GST}", 2)
$d.Content.Find.Execute("{#items}{quantity}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{quantity}", 2)
$d.Content.Find.Execute("{#items}{price}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{price}", 2)
$d.Content.Find.Execute("{#items}{finalAmount}{/items}", $true, $false, $false, $false, $false, $true, 1, $false, "{finalAmount}{/items}", 2)
